# ValueSet-uscore-treatment-intervention-preference.xlsx
# - Update a few Metadata values (ValueSet URL, Date)
# - Rename the existing "Include from US Core Treatmen" sheet to
#   "Include ValueSets" and reduce it to the ValueSet-URL row pair.
# - Add a brand-new "Include from US Core Treatmen" sheet after it,
#   carrying the concept/system-URI table (with the System URI value
#   repointed at hl7.org).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: update the ValueSet URL and the Date value.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item(1)

$meta.Range("B2").Value = "http://hl7.org/fhir/us/core/ValueSet/uscore-treatment-intervention-preference"

# Date needs to stay a plain text shared-string ("2023-10-02"), not get
# auto-converted into a date serial by Excel's type inference, so force
# text format first, then restore the untouched style (style index 2,
# as used by the neighbouring B7 cell) via a formats-only paste.
$meta.Range("B8").NumberFormat = "@"
$meta.Range("B8").Value = "2023-10-02"
$meta.Range("B7").Copy()
$meta.Range("B8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Rename the current second sheet to "Include ValueSets" and trim it
#    down to just the ValueSet URL row pair.
# ---------------------------------------------------------------------
$includeValueSets = $wb.Worksheets.Item(2)
$includeValueSets.Name = "Include ValueSets"

$includeValueSets.Range("A1").Value = "ValueSet URL"
$includeValueSets.Range("A2").Value = "https://cts.nlm.nih.gov/fhir/res/ValueSet/2.16.840.1.113762.1.4.1115.9"
$includeValueSets.Range("A3:B4").Clear()

# ---------------------------------------------------------------------
# 3. Add the new "Include from US Core Treatmen" sheet after it, with
#    the concept / description / system-URI table.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$includeCodes = $wb.Worksheets.Add($null, $lastSheet)
$includeCodes.Name = "Include from US Core Treatmen"

$includeCodes.Columns.Item(1).ColumnWidth = 29.8
$includeCodes.Columns.Item(2).ColumnWidth = 49.8

$includeCodes.Range("A1").Value = "Concept"
$includeCodes.Range("B1").Value = "Description"
$includeCodes.Range("A2").Value = "intervention-preference"
$includeCodes.Range("B2").Value = "Intervention preference"
$includeCodes.Range("A3").Value = ""
$includeCodes.Range("B3").Value = ""
$includeCodes.Range("A4").Value = "System URI"
$includeCodes.Range("B4").Value = "http://hl7.org/fhir/us/core/CodeSystem/uscore-treatment-intervention-preference"

# Match formatting: header row bold/shaded (style 1), body rows (style 2).
$includeValueSets.Range("A1").Copy()
$includeCodes.Range("A1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$includeCodes.Range("B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$includeValueSets.Range("A2").Copy()
$includeCodes.Range("A2:B4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Restore the active sheet/tab back to Metadata.
# ---------------------------------------------------------------------
$meta.Activate()
$meta.Range("A1").Select()
